$wb = $excel.ActiveWorkbook

# ===== Sheet ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 203.5
$ws.Range("I9").Value = 157
$ws.Range("K9").Value = 157
$ws.Range("M9").Value = 12
# Row 17
$ws.Range("H17").Value = 2039.6
$ws.Range("J17").Value = 2039.6
$ws.Range("L17").Value = 6118.799999999999
$ws.Range("N17").Value = -6454.799999999999
# Row 18
$ws.Range("H18").Value = 76925310
$ws.Range("I18").Value = 2069
$ws.Range("K18").Value = 2069
$ws.Range("M18").Value = -1785
# Row 28
$ws.Range("H28").Value = 1246.8334
$ws.Range("I28").Value = 435.3125
$ws.Range("J28").Value = 2869.875
$ws.Range("K28").Value = 435.3125
$ws.Range("L28").Value = 2869.875
$ws.Range("M28").Value = 49.6875
$ws.Range("N28").Value = -3839.875
# Row 62
$ws.Range("H62").Value = 2604.4119
$ws.Range("I62").Value = 2654.75
$ws.Range("J62").Value = 1799
$ws.Range("K62").Value = 2654.75
$ws.Range("L62").Value = 1799
$ws.Range("M62").Value = -2030.75
$ws.Range("N62").Value = -3047
# Row 64
$ws.Range("H64").Value = 142864420
$ws.Range("I64").Value = 8483.333000000001
$ws.Range("J64").Value = 1000000000
$ws.Range("K64").Value = 8483.333000000001
$ws.Range("L64").Value = 1000000000
$ws.Range("M64").Value = -8235.333000000001
$ws.Range("N64").Value = -1000000496
# Row 65
$ws.Range("H65").Value = 2604.4119
$ws.Range("I65").Value = 2654.75
$ws.Range("J65").Value = 1799
$ws.Range("K65").Value = 13273.75
$ws.Range("L65").Value = 8995
$ws.Range("M65").Value = -10153.75
$ws.Range("N65").Value = -15235
# Row 67
$ws.Range("H67").Value = 142864420
$ws.Range("I67").Value = 8483.333000000001
$ws.Range("J67").Value = 1000000000
$ws.Range("K67").Value = 8483.333000000001
$ws.Range("L67").Value = 1000000000
$ws.Range("M67").Value = -7625.333000000001
$ws.Range("N67").Value = -1000001716
# Row 98
$ws.Range("H98").Value = 1908.0889
$ws.Range("I98").Value = 1913.1395
$ws.Range("K98").Value = 1913.1395
$ws.Range("M98").Value = -415.1395
# Row 113
$ws.Range("H113").Value = 2614.4546
$ws.Range("I113").Value = 1961.8
$ws.Range("J113").Value = 3158.3333
$ws.Range("K113").Value = 1961.8
$ws.Range("L113").Value = 3158.3333
$ws.Range("M113").Value = 1292.2
$ws.Range("N113").Value = -9666.3333
# Row 121
$ws.Range("H121").Value = 2000
$ws.Range("J121").Value = 2000
$ws.Range("L121").Value = 6000
$ws.Range("N121").Value = -9494
# Row 122
$ws.Range("H122").Value = 1908.0889
$ws.Range("I122").Value = 1913.1395
$ws.Range("K122").Value = 5739.4185
$ws.Range("M122").Value = -3289.4185
# Row 132
$ws.Range("H132").Value = 7298.4546
$ws.Range("I132").Value = 7298.4546
$ws.Range("K132").Value = 21895.3638
$ws.Range("M132").Value = -19365.3638
# Row 138
$ws.Range("H138").Value = 332166.44
$ws.Range("I138").Value = 3482.111
$ws.Range("J138").Value = 513278.22
$ws.Range("K138").Value = 10446.333
$ws.Range("L138").Value = 1539834.66
$ws.Range("M138").Value = -5306.332999999999
$ws.Range("N138").Value = -1550114.66

# ===== Sheet ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 5154.636
$ws.Range("I61").Value = 3114.3333
$ws.Range("J61").Value = 7603
$ws.Range("K61").Value = 3114.3333
$ws.Range("L61").Value = 7603
$ws.Range("M61").Value = -2902.3333
$ws.Range("N61").Value = -8027
# Row 74
$ws.Range("H74").Value = 186753.5
$ws.Range("I74").Value = 278930.34
$ws.Range("J74").Value = 2399.8
$ws.Range("K74").Value = 278930.34
$ws.Range("L74").Value = 2399.8
$ws.Range("M74").Value = -278056.34
$ws.Range("N74").Value = -4147.8
# Row 77
$ws.Range("H77").Value = 186753.5
$ws.Range("I77").Value = 278930.34
$ws.Range("J77").Value = 2399.8
$ws.Range("K77").Value = 1394651.7
$ws.Range("L77").Value = 11999
$ws.Range("M77").Value = -1390283.7
$ws.Range("N77").Value = -20735
# Row 88
$ws.Range("H88").Value = 3451.0908
$ws.Range("J88").Value = 6969.75
$ws.Range("L88").Value = 6969.75
$ws.Range("N88").Value = -7781.75
# Row 91
$ws.Range("H91").Value = 3451.0908
$ws.Range("J91").Value = 6969.75
$ws.Range("L91").Value = 6969.75
$ws.Range("N91").Value = -9777.75
# Row 102
$ws.Range("H102").Value = 3183.2703
$ws.Range("I102").Value = 2566.1667
$ws.Range("K102").Value = 2566.1667
$ws.Range("M102").Value = -944.1667000000002
# Row 122
$ws.Range("H122").Value = 3917.6155
$ws.Range("I122").Value = 3613.9395
$ws.Range("K122").Value = 10841.8185
$ws.Range("M122").Value = -8391.818499999999
# Row 132
$ws.Range("H132").Value = 3040.8286
$ws.Range("I132").Value = 1893
$ws.Range("J132").Value = 5545.1816
$ws.Range("K132").Value = 5679
$ws.Range("L132").Value = 16635.5448
$ws.Range("M132").Value = -3149
$ws.Range("N132").Value = -21695.5448
# Row 136
$ws.Range("H136").Value = 5154.636
$ws.Range("I136").Value = 3114.3333
$ws.Range("J136").Value = 7603
$ws.Range("K136").Value = 9342.999899999999
$ws.Range("L136").Value = 22809
$ws.Range("M136").Value = -6792.999899999999
$ws.Range("N136").Value = -27909

# ===== Sheet BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 737.4286
$ws.Range("J64").Value = 1092.3334
$ws.Range("L64").Value = 1092.3334
$ws.Range("N64").Value = -1542.3334
# Row 67
$ws.Range("H67").Value = 737.4286
$ws.Range("J67").Value = 1092.3334
$ws.Range("L67").Value = 1092.3334
$ws.Range("N67").Value = -2652.3334
# Row 94
$ws.Range("H94").Value = 153846460
$ws.Range("I94").Value = 153846460
$ws.Range("K94").Value = 153846460
$ws.Range("M94").Value = -153846009
# Row 134
$ws.Range("H134").Value = 2751.1072
$ws.Range("I134").Value = 2384.05
$ws.Range("J134").Value = 3668.75
$ws.Range("K134").Value = 7152.150000000001
$ws.Range("L134").Value = 11006.25
$ws.Range("M134").Value = -4617.150000000001
$ws.Range("N134").Value = -16076.25

# ===== Sheet CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 4831.769
$ws.Range("I99").Value = 4181.3
$ws.Range("J99").Value = 7000
$ws.Range("K99").Value = 4181.3
$ws.Range("L99").Value = 7000
$ws.Range("M99").Value = -2683.3
$ws.Range("N99").Value = -9996
# Row 126
$ws.Range("H126").Value = 4831.769
$ws.Range("I126").Value = 4181.3
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 12543.9
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -10073.9
$ws.Range("N126").Value = -25940

# ===== Sheet CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1976.9546
$ws.Range("I5").Value = 428
$ws.Range("J5").Value = 2221.5264
$ws.Range("K5").Value = 1284
$ws.Range("L5").Value = 6664.5792
$ws.Range("M5").Value = -1172
$ws.Range("N5").Value = -6888.5792
# Row 14
$ws.Range("H14").Value = 564.63635
$ws.Range("I14").Value = 564.63635
$ws.Range("K14").Value = 1693.90905
$ws.Range("M14").Value = -1520.90905
# Row 97
$ws.Range("H97").Value = 833763.3
$ws.Range("J97").Value = 525
$ws.Range("L97").Value = 1575
$ws.Range("N97").Value = -2567
# Row 107
$ws.Range("H107").Value = 529.8889
$ws.Range("I107").Value = 322
$ws.Range("J107").Value = 589.2857
$ws.Range("K107").Value = 966
$ws.Range("L107").Value = 1767.8571
$ws.Range("M107").Value = 954
$ws.Range("N107").Value = -5607.8571
# Row 129
$ws.Range("H129").Value = 1369
$ws.Range("I129").Value = 939
$ws.Range("K129").Value = 2817
$ws.Range("M129").Value = 2183
# Row 135
$ws.Range("H135").Value = 1976.9546
$ws.Range("I135").Value = 428
$ws.Range("J135").Value = 2221.5264
$ws.Range("K135").Value = 3852
$ws.Range("L135").Value = 19993.7376
$ws.Range("M135").Value = -1317
$ws.Range("N135").Value = -25063.7376

# ===== Sheet GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 200003100
$ws.Range("I80").Value = 1000000000
$ws.Range("J80").Value = 3874.5
$ws.Range("K80").Value = 1000000000
$ws.Range("L80").Value = 3874.5
$ws.Range("M80").Value = -999999002
$ws.Range("N80").Value = -5870.5
# Row 83
$ws.Range("H83").Value = 200003100
$ws.Range("I83").Value = 1000000000
$ws.Range("J83").Value = 3874.5
$ws.Range("K83").Value = 5000000000
$ws.Range("L83").Value = 19372.5
$ws.Range("M83").Value = -4999995008
$ws.Range("N83").Value = -29356.5
# Row 126
$ws.Range("H126").Value = 3577.4092
$ws.Range("I126").Value = 1338.4667
$ws.Range("J126").Value = 8375.143
$ws.Range("K126").Value = 4015.4001
$ws.Range("L126").Value = 25125.429
$ws.Range("M126").Value = -1545.4001
$ws.Range("N126").Value = -30065.429

# ===== Sheet LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1708.6666
$ws.Range("I46").Value = 750
$ws.Range("J46").Value = 2188
$ws.Range("K46").Value = 750
$ws.Range("L46").Value = 2188
$ws.Range("M46").Value = -562
$ws.Range("N46").Value = -2564
# Row 61
$ws.Range("H61").Value = 3401
$ws.Range("I61").Value = 3401
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3401
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3199
$ws.Range("N61").ClearContents()
# Row 113
$ws.Range("H113").Value = 3401
$ws.Range("I113").Value = 3401
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3401
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1231
$ws.Range("N113").ClearContents()
# Row 121
$ws.Range("H121").Value = 54420
$ws.Range("J121").Value = 54420
$ws.Range("L121").Value = 54420
$ws.Range("N121").Value = -57914
# Row 122
$ws.Range("H122").Value = 3304.4443
$ws.Range("I122").Value = 3280
$ws.Range("K122").Value = 9840
$ws.Range("M122").Value = -7390

# ===== Sheet WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 1348.5555
$ws.Range("I107").Value = 876.8570999999999
$ws.Range("J107").Value = 2999.5
$ws.Range("K107").Value = 2630.5713
$ws.Range("L107").Value = 8998.5
$ws.Range("M107").Value = -710.5712999999996
$ws.Range("N107").Value = -12838.5
# Row 121
$ws.Range("H121").Value = 110000
$ws.Range("J121").Value = 110000
$ws.Range("L121").Value = 110000
$ws.Range("N121").Value = -113494
